# Apply odds/stat updates for 2025-05-06 FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("P13").Value = 1.7
$ws.Range("W13").Value = 25
$ws.Range("AE13").Value = 6.4

# Row 27
$ws.Range("G27").Value = 2.9
$ws.Range("I27").Value = 2.4
$ws.Range("J27").Value = 1.08
$ws.Range("K27").Value = 8
$ws.Range("T27").Value = 8
$ws.Range("U27").Value = 13
$ws.Range("V27").Value = 11
$ws.Range("W27").Value = 29
$ws.Range("AF27").Value = 11
$ws.Range("AH27").Value = 23

# Row 41
$ws.Range("G41").Value = 4.33
$ws.Range("H41").Value = 3.3
$ws.Range("I41").Value = 1.9
$ws.Range("T41").Value = 10
$ws.Range("AC41").Value = 51

# Row 42
$ws.Range("N42").Value = 1.9
$ws.Range("O42").Value = 1.95
